# Insert a new weekly data point at row 191 (a new Coliflor price record for
# Terminal Hortofrutícola Agro Chillán), pushing the existing rows 191-245
# down to 192-246. This mirrors how the source publishes data with the
# newest observation inserted near the top of the date-ordered block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 191:245 down to 192:246, leaving a blank row 191 to populate.
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Cells.Item(191, 1).Value  = 7
$ws.Cells.Item(191, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(191, 3).Value  = "Ñuble"
$ws.Cells.Item(191, 4).Value  = 44642
$ws.Cells.Item(191, 5).Value  = 16
$ws.Cells.Item(191, 6).Value  = 100112008
$ws.Cells.Item(191, 7).Value  = "Coliflor"
$ws.Cells.Item(191, 8).Value  = "Sin especificar"
$ws.Cells.Item(191, 9).Value  = "Primera"
$ws.Cells.Item(191, 10).Value = 200
$ws.Cells.Item(191, 11).Value = 1300
$ws.Cells.Item(191, 12).Value = 1400
$ws.Cells.Item(191, 13).Value = 1350
$ws.Cells.Item(191, 14).Value = "`$/unidad"
$ws.Cells.Item(191, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(191, 16).Value = 1350
$ws.Cells.Item(191, 17).Value = 1
$ws.Cells.Item(191, 18).Value = "Hortaliza"
